$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.506.28"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.540.16"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.46"
$ws.Range("E5").Value = "  +4.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.10"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "3.535.17"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.199"
$ws.Range("E10").Value = "  +5.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.78"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.585"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.27"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "4.119.20"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "628.99"
$ws.Range("E16").Value = "  -6.63%  "
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.543.14"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "70.509.23"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.41"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.01"
$ws.Range("E22").Value = "  -10.92%  "
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.93"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.69"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.62"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.44"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.49"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.98"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "570.65"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.62"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.49"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0464"
$ws.Range("E41").Value = "  +6.47%  "
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").Value = "3.338.40"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("E45").Value = "  +5.56%  "
$ws.Range("D46").Value = "0.0₃0713"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.08"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.64"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.70"
$ws.Range("E51").Value = "  +1.77%  "
